$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.282184362411499
$ws.Range("B1").Value = 2.432773590087891
$ws.Range("C1").Value = 3.306196451187134
$ws.Range("D1").Value = 3.27686595916748
$ws.Range("E1").Value = 1.065208554267883
